# Apply the "Updated ITA model - 2025-09-01 00:07" edit to Scen_Base_VS.xlsx
#
# Summary of the change (per the OOXML diff):
#  - Veda!B6  "solar" -> "\I: solar"   (VLOOKUP in D6 now misses -> D6 becomes "")
#  - Veda!B7  "wind"  -> "\I: wind"    (VLOOKUP in D7 now misses -> D7 becomes "")
#  - Veda!F11 "wind"  -> "windon"      (renaming the wind-onshore row of the lookup table)
#  - historical_data_long!A11,A20,A29,...,A650 (every 9th row) "wind" -> "windon"
#  - The active sheet/selection moves from buildrates!D15 to Veda!B8
#    (previously buildrates was the active tab with Veda!V3 selected)

$wb = $excel.ActiveWorkbook

$wsVeda = $wb.Worksheets.Item("Veda")
$wsBuildrates = $wb.Worksheets.Item("buildrates")
$wsHist = $wb.Worksheets.Item("historical_data_long")

# --- Veda sheet: rename lookup keys -------------------------------------------------
$wsVeda.Range("B6").Value = "\I: solar"
$wsVeda.Range("B7").Value = "\I: wind"
$wsVeda.Range("F11").Value = "windon"

# --- historical_data_long sheet: rename every "wind" row label to "windon" ----------
$histRow = 11
while ($histRow -le 650) {
    $wsHist.Range("A$histRow").Value = "windon"
    $histRow = $histRow + 9
}

# Recalculate so the cached formula results (e.g. Veda!D6 / Veda!D7) reflect the
# changed lookup keys.
$excel.Calculate()

# --- Selection / active sheet changes ------------------------------------------------
# Previously buildrates was the active tab (selection D15, which stays unchanged);
# now Veda is the active tab with B8 selected.
$wsVeda.Activate()
$wsVeda.Range("B8").Select()
